$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("multicolsWithNulls")
$ws.Activate()

# Scratch area (well outside the used range) used to stage a copy of one
# side of the swap so it isn't lost while the destination is overwritten.
# Copy + PasteSpecial(xlPasteValues = -4163) transfers the literal value
# (and its type: number vs text) of each cell without disturbing the
# destination cell's own formatting/style index.
$xlPasteValues = -4163

function Swap-Range([string]$rangeA, [string]$rangeB, [string]$scratch) {
    $ws.Range($rangeA).Copy($ws.Range($scratch))
    $ws.Range($rangeB).Copy()
    $ws.Range($rangeA).PasteSpecial($xlPasteValues)
    $ws.Range($scratch).Copy()
    $ws.Range($rangeB).PasteSpecial($xlPasteValues)
    $ws.Range($scratch).Clear()
    $excel.CutCopyMode = 0
}

# ---------------------------------------------------------------------
# Row 5 <-> Row 6 swap (columns B, D, E, F, G, H, I).
# Column C is identical in both rows ("12") so it is left untouched.
# ---------------------------------------------------------------------
Swap-Range "B5" "B6" "K1"
Swap-Range "D5:F5" "D6:F6" "K1:M1"
Swap-Range "G5:I5" "G6:I6" "K1:M1"

# ---------------------------------------------------------------------
# Row 9 <-> Row 10 swap (columns C, D, E, F, G, H, I).
# Column B is identical in both rows ("11.0") so it is left untouched.
# ---------------------------------------------------------------------
Swap-Range "C9" "C10" "K1"
Swap-Range "D9:F9" "D10:F10" "K1:M1"
Swap-Range "G9:I9" "G10:I10" "K1:M1"
